# IKD update: GaN CMOS 2026-01-31T23:27Z
# Append 6 new literature records (rows 113-118) to the Master sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-LitRow($Row, $Title, $Year, $Publisher, $Venue, $Authors, $Organizations, $DOI, $Url, $DocType, $DeviceType, $Method, $EnablerCategory, $MaterialSystem, $NodeGeometry, $KeyContribution, $EvidenceSnippet, $TagConfidence, $AddedDate, $Notes) {
    $ws.Cells.Item($Row, 1).Value = ""
    $ws.Cells.Item($Row, 2).Value = $Title
    $ws.Cells.Item($Row, 3).Value = $Year
    $ws.Cells.Item($Row, 4).Value = $Publisher
    $ws.Cells.Item($Row, 5).Value = $Venue
    $ws.Cells.Item($Row, 6).Value = $Authors
    $ws.Cells.Item($Row, 7).Value = $Organizations
    $ws.Cells.Item($Row, 8).Value = $DOI
    $ws.Cells.Item($Row, 9).Value = $Url
    $ws.Cells.Item($Row, 10).Value = $DocType
    $ws.Cells.Item($Row, 11).Value = $DeviceType
    $ws.Cells.Item($Row, 12).Value = $Method
    $ws.Cells.Item($Row, 13).Value = $EnablerCategory
    $ws.Cells.Item($Row, 14).Value = $MaterialSystem
    $ws.Cells.Item($Row, 15).Value = $NodeGeometry
    $ws.Cells.Item($Row, 16).Value = $KeyContribution
    $ws.Cells.Item($Row, 17).Value = $EvidenceSnippet
    $ws.Cells.Item($Row, 18).Value = $TagConfidence
    # Force AddedDate to remain plain text (matches existing rows), not an
    # auto-converted date serial value.
    $ws.Cells.Item($Row, 19).NumberFormat = "@"
    $ws.Cells.Item($Row, 19).Value = $AddedDate
    $ws.Cells.Item($Row, 20).Value = ""
}

# Row 113 - CMOS Measurement-Collapse Primitive (1 of 2)
Set-LitRow 113 `
    "A CMOS Measurement-Collapse Primitive for Ephemeral Secrets in Post-Quantum Cryptography" `
    2026 `
    "Institute of Electrical and Electronics Engineers (IEEE)" `
    "" `
    "III, Francis X. Cunnane" `
    "" `
    "10.36227/techrxiv.176463742.23048082/v3" `
    "https://doi.org/10.36227/techrxiv.176463742.23048082/v3" `
    "Journal" `
    "Co-integration" `
    "Experiment" `
    "Gate Stack" `
    "" `
    "" `
    "" `
    "A CMOS Measurement-Collapse Primitive for Ephemeral Secrets in Post-Quantum Cryptography" `
    "High" `
    "2026-01-31" `
    ""

# Row 114 - CMOS Measurement-Collapse Primitive (2 of 2, duplicate)
Set-LitRow 114 `
    "A CMOS Measurement-Collapse Primitive for Ephemeral Secrets in Post-Quantum Cryptography" `
    2026 `
    "Institute of Electrical and Electronics Engineers (IEEE)" `
    "" `
    "III, Francis X. Cunnane" `
    "" `
    "10.36227/techrxiv.176463742.23048082/v3" `
    "https://doi.org/10.36227/techrxiv.176463742.23048082/v3" `
    "Journal" `
    "Co-integration" `
    "Experiment" `
    "Gate Stack" `
    "" `
    "" `
    "" `
    "A CMOS Measurement-Collapse Primitive for Ephemeral Secrets in Post-Quantum Cryptography" `
    "High" `
    "2026-01-31" `
    ""

# Row 115 - TCAD Ga2O3 MISFET (1 of 4)
Set-LitRow 115 `
    "TCAD Demonstration of a High-Voltage Lateral Double-RESURF-VLD β-Ga2O3 MISFET with p-type Diamond for Ultra-Low Ron,sp" `
    2026 `
    "The Electrochemical Society" `
    "ECS Journal of Solid State Science and Technology" `
    "Xu, Fan; Dai, Ming; Duan, JunFeng; Zhu, Shengnan; Qiao, Yuan; Yi, Bo; Cheng, Junji; Yang, Hongqiang" `
    "" `
    "10.1149/2162-8777/ae3fce" `
    "https://doi.org/10.1149/2162-8777/ae3fce" `
    "Journal" `
    "n-FET" `
    "TCAD" `
    "Gate Stack" `
    "" `
    "" `
    "" `
    "TCAD Demonstration of a High-Voltage Lateral Double-RESURF-VLD β-Ga2O3 MISFET with p-type Diamond for Ultra-Low Ron,sp" `
    "High" `
    "2026-01-31" `
    ""

# Row 116 - TCAD Ga2O3 MISFET (2 of 4, duplicate)
Set-LitRow 116 `
    "TCAD Demonstration of a High-Voltage Lateral Double-RESURF-VLD β-Ga2O3 MISFET with p-type Diamond for Ultra-Low Ron,sp" `
    2026 `
    "The Electrochemical Society" `
    "ECS Journal of Solid State Science and Technology" `
    "Xu, Fan; Dai, Ming; Duan, JunFeng; Zhu, Shengnan; Qiao, Yuan; Yi, Bo; Cheng, Junji; Yang, Hongqiang" `
    "" `
    "10.1149/2162-8777/ae3fce" `
    "https://doi.org/10.1149/2162-8777/ae3fce" `
    "Journal" `
    "n-FET" `
    "TCAD" `
    "Gate Stack" `
    "" `
    "" `
    "" `
    "TCAD Demonstration of a High-Voltage Lateral Double-RESURF-VLD β-Ga2O3 MISFET with p-type Diamond for Ultra-Low Ron,sp" `
    "High" `
    "2026-01-31" `
    ""

# Row 117 - TCAD Ga2O3 MISFET (3 of 4, duplicate)
Set-LitRow 117 `
    "TCAD Demonstration of a High-Voltage Lateral Double-RESURF-VLD β-Ga2O3 MISFET with p-type Diamond for Ultra-Low Ron,sp" `
    2026 `
    "The Electrochemical Society" `
    "ECS Journal of Solid State Science and Technology" `
    "Xu, Fan; Dai, Ming; Duan, JunFeng; Zhu, Shengnan; Qiao, Yuan; Yi, Bo; Cheng, Junji; Yang, Hongqiang" `
    "" `
    "10.1149/2162-8777/ae3fce" `
    "https://doi.org/10.1149/2162-8777/ae3fce" `
    "Journal" `
    "n-FET" `
    "TCAD" `
    "Gate Stack" `
    "" `
    "" `
    "" `
    "TCAD Demonstration of a High-Voltage Lateral Double-RESURF-VLD β-Ga2O3 MISFET with p-type Diamond for Ultra-Low Ron,sp" `
    "High" `
    "2026-01-31" `
    ""

# Row 118 - TCAD Ga2O3 MISFET (4 of 4, duplicate)
Set-LitRow 118 `
    "TCAD Demonstration of a High-Voltage Lateral Double-RESURF-VLD β-Ga2O3 MISFET with p-type Diamond for Ultra-Low Ron,sp" `
    2026 `
    "The Electrochemical Society" `
    "ECS Journal of Solid State Science and Technology" `
    "Xu, Fan; Dai, Ming; Duan, JunFeng; Zhu, Shengnan; Qiao, Yuan; Yi, Bo; Cheng, Junji; Yang, Hongqiang" `
    "" `
    "10.1149/2162-8777/ae3fce" `
    "https://doi.org/10.1149/2162-8777/ae3fce" `
    "Journal" `
    "n-FET" `
    "TCAD" `
    "Gate Stack" `
    "" `
    "" `
    "" `
    "TCAD Demonstration of a High-Voltage Lateral Double-RESURF-VLD β-Ga2O3 MISFET with p-type Diamond for Ultra-Low Ron,sp" `
    "High" `
    "2026-01-31" `
    ""
